$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.885083913803101
$ws.Range("B1").Value = 5.206869602203369
$ws.Range("C1").Value = 6.728024482727051
$ws.Range("D1").Value = 11.30191516876221
$ws.Range("E1").Value = 4.450334548950195
